$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ANTHONY_MBM_Worked")

# New ticket-assignment event appended to the worked log
$ws.Range("A56").Value = "17:27:06 04-12-2023"
$ws.Range("B56").Value = "Automatically Assigned Ticket"

# Widen column A so the longer timestamp strings stay readable
$ws.Columns.Item(1).ColumnWidth = 20.592447916666668

# Bring the new row into view / select it, like a user scrolling to it
[void]$ws.Range("A57:E57").Select()
